$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 560.74536596201631
$ws.Range("C2").Value = 398.61166425469929
$ws.Range("D2").Value = 715.43504452547893
$ws.Range("E2").Value = 393.25617965248244

$ws.Range("B3").Value = 702.03251463894992
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = 779.2615068409574
$ws.Range("E3").Value = 402.15337977147959

$ws.Range("B1:E3").Select()
